$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.956950000000001
$ws.Range("H2").Value = 26.87085
$ws.Range("I2").Value = 0.3465211830970586
$ws.Range("J2").Value = 0.3465211830970586
$ws.Range("M2").Value = 1021.934916333333
$ws.Range("N2").Value = 3065.804749
$ws.Range("O2").Value = 0.8026347959846111
$ws.Range("P2").Value = 0.802634795984611
$ws.Range("Q2").Value = 9153.41994885185
$ws.Range("R2").Value = 82380.77953966665
$ws.Range("S2").Value = 0.2781299590994537
$ws.Range("T2").Value = 0.2781299590994536
$ws.Range("G3").Value = 8.956950000000001
$ws.Range("H3").Value = 26.87085
$ws.Range("I3").Value = 0.3465211830970586
$ws.Range("J3").Value = 0.3465211830970586
$ws.Range("O3").Value = 0.04931810976893385
$ws.Range("P3").Value = 0.04931810976893384
$ws.Range("Q3").Value = 562.4343375804501
$ws.Range("R3").Value = 5061.90903822405
$ws.Range("S3").Value = 0.01708976974524156
$ws.Range("T3").Value = 0.01708976974524156
$ws.Range("G4").Value = 8.956950000000001
$ws.Range("H4").Value = 26.87085
$ws.Range("I4").Value = 0.3465211830970586
$ws.Range("J4").Value = 0.3465211830970586
$ws.Range("M4").Value = 187.139577
$ws.Range("N4").Value = 561.418731
$ws.Range("O4").Value = 0.1469807262726385
$ws.Range("P4").Value = 0.1469807262726385
$ws.Range("Q4").Value = 1676.19983421015
$ws.Range("R4").Value = 15085.79850789135
$ws.Range("S4").Value = 0.05093193516045963
$ws.Range("T4").Value = 0.05093193516045962
$ws.Range("G5").Value = 8.956950000000001
$ws.Range("H5").Value = 26.87085
$ws.Range("I5").Value = 0.3465211830970586
$ws.Range("J5").Value = 0.3465211830970586
$ws.Range("M5").Value = 1.357726666666667
$ws.Range("N5").Value = 4.073180000000001
$ws.Range("O5").Value = 0.001066367973816652
$ws.Range("P5").Value = 0.001066367973816652
$ws.Range("Q5").Value = 12.161089867
$ws.Range("R5").Value = 109.449808803
$ws.Range("S5").Value = 0.0003695190919037595
$ws.Range("T5").Value = 0.0003695190919037593
$ws.Range("I6").Value = 0.2466462208011621
$ws.Range("J6").Value = 0.2466462208011621
$ws.Range("M6").Value = 1021.934916333333
$ws.Range("N6").Value = 3065.804749
$ws.Range("O6").Value = 0.8026347959846111
$ws.Range("P6").Value = 0.802634795984611
$ws.Range("Q6").Value = 6515.20469141974
$ws.Range("R6").Value = 58636.84222277766
$ws.Range("S6").Value = 0.197966839113116
$ws.Range("T6").Value = 0.197966839113116
$ws.Range("I7").Value = 0.2466462208011621
$ws.Range("J7").Value = 0.2466462208011621
$ws.Range("O7").Value = 0.04931810976893385
$ws.Range("P7").Value = 0.04931810976893384
$ws.Range("S7").Value = 0.01216412539156441
$ws.Range("T7").Value = 0.0121641253915644
$ws.Range("I8").Value = 0.2466462208011621
$ws.Range("J8").Value = 0.2466462208011621
$ws.Range("M8").Value = 187.139577
$ws.Range("N8").Value = 561.418731
$ws.Range("O8").Value = 0.1469807262726385
$ws.Range("P8").Value = 0.1469807262726385
$ws.Range("Q8").Value = 1193.082485522015
$ws.Range("R8").Value = 10737.74236969814
$ws.Range("S8").Value = 0.03625224066575636
$ws.Range("T8").Value = 0.03625224066575636
$ws.Range("I9").Value = 0.2466462208011621
$ws.Range("J9").Value = 0.2466462208011621
$ws.Range("M9").Value = 1.357726666666667
$ws.Range("N9").Value = 4.073180000000001
$ws.Range("O9").Value = 0.001066367973816652
$ws.Range("P9").Value = 0.001066367973816652
$ws.Range("Q9").Value = 8.65599854447778
$ws.Range("R9").Value = 77.90398690030001
$ws.Range("S9").Value = 0.0002630156307252698
$ws.Range("T9").Value = 0.0002630156307252697
$ws.Range("G10").Value = 7.505276333333332
$ws.Range("H10").Value = 22.515829
$ws.Range("I10").Value = 0.2903596910217228
$ws.Range("J10").Value = 0.2903596910217228
$ws.Range("M10").Value = 1021.934916333333
$ws.Range("N10").Value = 3065.804749
$ws.Range("O10").Value = 0.8026347959846111
$ws.Range("P10").Value = 0.802634795984611
$ws.Range("Q10").Value = 7669.903941763545
$ws.Range("R10").Value = 69029.13547587191
$ws.Range("S10").Value = 0.2330527913653752
$ws.Range("T10").Value = 0.2330527913653751
$ws.Range("G11").Value = 7.505276333333332
$ws.Range("H11").Value = 22.515829
$ws.Range("I11").Value = 0.2903596910217228
$ws.Range("J11").Value = 0.2903596910217228
$ws.Range("O11").Value = 0.04931810976893385
$ws.Range("P11").Value = 0.04931810976893384
$ws.Range("Q11").Value = 471.2792996384441
$ws.Range("R11").Value = 4241.513696745997
$ws.Range("S11").Value = 0.01431999111428304
$ws.Range("T11").Value = 0.01431999111428304
$ws.Range("G12").Value = 7.505276333333332
$ws.Range("H12").Value = 22.515829
$ws.Range("I12").Value = 0.2903596910217228
$ws.Range("J12").Value = 0.2903596910217228
$ws.Range("M12").Value = 187.139577
$ws.Range("N12").Value = 561.418731
$ws.Range("O12").Value = 0.1469807262726385
$ws.Range("P12").Value = 0.1469807262726385
$ws.Range("Q12").Value = 1404.534238288111
$ws.Range("R12").Value = 12640.808144593
$ws.Range("S12").Value = 0.04267727826667173
$ws.Range("T12").Value = 0.04267727826667173
$ws.Range("G13").Value = 7.505276333333332
$ws.Range("H13").Value = 22.515829
$ws.Range("I13").Value = 0.2903596910217228
$ws.Range("J13").Value = 0.2903596910217228
$ws.Range("M13").Value = 1.357726666666667
$ws.Range("N13").Value = 4.073180000000001
$ws.Range("O13").Value = 0.001066367973816652
$ws.Range("P13").Value = 0.001066367973816652
$ws.Range("Q13").Value = 10.19011381846889
$ws.Range("R13").Value = 91.71102436622
$ws.Range("S13").Value = 0.0003096302753928637
$ws.Range("T13").Value = 0.0003096302753928635
$ws.Range("G14").Value = 3.010615333333333
$ws.Range("H14").Value = 9.031846
$ws.Range("I14").Value = 0.1164729050800565
$ws.Range("J14").Value = 0.1164729050800565
$ws.Range("M14").Value = 1021.934916333333
$ws.Range("N14").Value = 3065.804749
$ws.Range("O14").Value = 0.8026347959846111
$ws.Range("P14").Value = 0.802634795984611
$ws.Range("Q14").Value = 3076.65292878185
$ws.Range("R14").Value = 27689.87635903665
$ws.Range("S14").Value = 0.09348520640666612
$ws.Range("T14").Value = 0.0934852064066661
$ws.Range("G15").Value = 3.010615333333333
$ws.Range("H15").Value = 9.031846
$ws.Range("I15").Value = 0.1164729050800565
$ws.Range("J15").Value = 0.1164729050800565
$ws.Range("O15").Value = 0.04931810976893385
$ws.Range("P15").Value = 0.04931810976893384
$ws.Range("Q15").Value = 189.0457623089198
$ws.Range("R15").Value = 1701.411860780278
$ws.Range("S15").Value = 0.00574422351784484
$ws.Range("T15").Value = 0.005744223517844838
$ws.Range("G16").Value = 3.010615333333333
$ws.Range("H16").Value = 9.031846
$ws.Range("I16").Value = 0.1164729050800565
$ws.Range("J16").Value = 0.1164729050800565
$ws.Range("M16").Value = 187.139577
$ws.Range("N16").Value = 561.418731
$ws.Range("O16").Value = 0.1469807262726385
$ws.Range("P16").Value = 0.1469807262726385
$ws.Range("Q16").Value = 563.405279989714
$ws.Range("R16").Value = 5070.647519907426
$ws.Range("S16").Value = 0.0171192721797508
$ws.Range("T16").Value = 0.01711927217975079
$ws.Range("G17").Value = 3.010615333333333
$ws.Range("H17").Value = 9.031846
$ws.Range("I17").Value = 0.1164729050800565
$ws.Range("J17").Value = 0.1164729050800565
$ws.Range("M17").Value = 1.357726666666667
$ws.Range("N17").Value = 4.073180000000001
$ws.Range("O17").Value = 0.001066367973816652
$ws.Range("P17").Value = 0.001066367973816652
$ws.Range("Q17").Value = 4.087592721142223
$ws.Range("R17").Value = 36.78833449028001
$ws.Range("S17").Value = 0.0001242029757947591
$ws.Range("T17").Value = 0.0001242029757947591

Write-Output "Applied 174 cell updates"